$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 updates
$ws.Range("D2").Value = 2649.026431718062
$ws.Range("E2").Value = 6938.179297474234
$ws.Range("J2").Value = 4438.114537444934
$ws.Range("K2").Value = 12094.84137329755
$ws.Range("P2").Value = 11266.90748898678
$ws.Range("Q2").Value = 25891.47595910595
$ws.Range("V2").Value = 5669.88986784141
$ws.Range("W2").Value = 12141.26430089142
$ws.Range("Y2").Value = 554
$ws.Range("Z2").Value = 5524.5
$ws.Range("AB2").Value = 4443.290748898678
$ws.Range("AC2").Value = 10451.28341757045

# Row 3 updates
$ws.Range("D3").Value = 8617.881057268722
$ws.Range("J3").Value = 6828.79295154185
$ws.Range("V3").Value = 5597.017621145375
$ws.Range("AB3").Value = 6823.616740088106
